$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56 (shifts existing rows 56..130 down to 57..131)
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly price record
# (2021-08-05, the day after the previous most-recent record).
$ws.Cells.Item(56, 1).Value = 9
$ws.Cells.Item(56, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(56, 3).Value = "Metropolitana"
$ws.Cells.Item(56, 4).Value = 44413
$ws.Cells.Item(56, 5).Value = 13
$ws.Cells.Item(56, 6).Value = 300000001
$ws.Cells.Item(56, 7).Value = "Rabanito"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 8800
$ws.Cells.Item(56, 11).Value = 2500
$ws.Cells.Item(56, 12).Value = 3000
$ws.Cells.Item(56, 13).Value = 2750
$ws.Cells.Item(56, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(56, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(56, 16).Value = 28
$ws.Cells.Item(56, 17).Value = 100
$ws.Cells.Item(56, 18).Value = "Hortaliza"
